# Remove the "surtitle" placeholder (a body placeholder at ph idx="10") that
# was defined on the "1_Surtitle, Title Only" slide layout and instantiated
# (with overridden text "Column_subtitle") on slide 2, which uses that layout.
#
# The edit deletes the shape both from the slide layout definition and from
# the slide's own override, leaving every other shape untouched.

$p = $ppt.ActivePresentation

# Slide 2 uses the "1_Surtitle, Title Only" layout and overrides the surtitle
# placeholder text with "Column_subtitle".
$s = $p.Slides.Item(2)
$layout = $s.CustomLayout

# --- 1) Remove the placeholder's definition from the slide layout ---
for ($i = $layout.Shapes.Count; $i -ge 1; $i--) {
    $shp = $layout.Shapes.Item($i)
    if ($shp.Name -eq "Text Placeholder 7") {
        $shp.Delete()
    }
}

# --- 2) Remove the corresponding placeholder instance from the slide ---
# The shape is an instantiated layout placeholder, so the first Delete() only
# clears its text/content: PowerPoint re-synthesizes an empty placeholder
# shape in its place for as long as the slide still references that
# placeholder idx. Deleting a second time removes it for good, matching the
# target (no surtitle shape left on the slide at all).
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "Text Placeholder 1") {
        $shp.Delete()
    }
}
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.PlaceholderFormat.Type -eq 2) {
        $shp.Delete()
    }
}
